$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Report")

# "Refactoring & Re-Work" row (row 23): add actual-cost entries for S/T/U (weeks 15-17)
$ws.Range("S23").Value = 3
$ws.Range("T23").Value = 3
$ws.Range("U23").Value = 3

# "Documentation" row (row 29): add actual-cost entries for S/T/U (weeks 15-17)
$ws.Range("S29").Value = 2
$ws.Range("T29").Value = 2
$ws.Range("U29").Value = 2
